$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "2022 Performance by Country" summary table (B2:I76) was copied and
# pasted as values into a hidden staging block (M3:T77, offset by +1 row /
# +11 columns) that feeds a chart.
$ws.Range("B2:I76").Copy()
$ws.Range("M3").PasteSpecial(-4163)

# Clear clipboard marching ants / pasted-range selection artifact.
$excel.CutCopyMode = $false

# Scroll/zoom state left behind by the author after finishing the paste.
$excel.ActiveWindow.Zoom = 46
$ws.Range("W25").Select()
